$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Best-effort workbook window geometry (cosmetic; may not persist, harmless).
# ---------------------------------------------------------------------------
$win = $wb.Windows.Item(1)
$win.Left = -13500
$win.Top = 3540
$win.Width = 13500
$win.Height = 14520

# ---------------------------------------------------------------------------
# Helper values shared by all three rubric sheets (rows 7-16, cols D:E).
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)   # "Group A"
$ws2 = $wb.Worksheets.Item(2)   # "Group B"
$ws3 = $wb.Worksheets.Item(3)   # "Group C"

foreach ($ws in @($ws1, $ws2, $ws3)) {
    # Part header text: row4 label changes, D4/E4 lose explicit style (keep the 10s)
    $ws.Range("D4:E4").ClearFormats()

    # Row5 + Row6 D/E cells were blank placeholders; fully clear them so the
    # <c> elements disappear from the XML entirely.
    $ws.Range("D5:E6").Clear()

    # Revised point totals for rows 7-16.
    $ws.Range("D7").Value = 2
    $ws.Range("E7").Value = 2
    $ws.Range("D8").Value = 3
    $ws.Range("E8").Value = 3
    $ws.Range("D9").Value = 4
    $ws.Range("E9").Value = 4
    $ws.Range("D10").Value = 4
    $ws.Range("E10").Value = 4
    $ws.Range("D11").Value = 4
    $ws.Range("E11").Value = 4
    $ws.Range("D12").Value = 2
    $ws.Range("E12").Value = 2
    $ws.Range("D13").Value = 4
    $ws.Range("E13").Value = 4
    $ws.Range("D15").Value = 2
    $ws.Range("E15").Value = 2
    $ws.Range("D16").Value = 2
    $ws.Range("E16").Value = 2
}

# ---------------------------------------------------------------------------
# Sheet-specific text relabels: "Part A: Exercise" -> "Part 1: Exercise" and
# "Part B" -> "Part 2" on Group A / Group C only (Group B keeps the old text).
# ---------------------------------------------------------------------------
$ws1.Range("A4").Value = "Part 1: Exercise"
$ws1.Range("A6").Value = "Part 2"

$ws3.Range("A4").Value = "Part 1: Exercise"
$ws3.Range("A6").Value = "Part 2"

# ---------------------------------------------------------------------------
# Column width tweaks.
# ---------------------------------------------------------------------------
$ws1.Columns.Item(3).ColumnWidth = 21.498697916666668
$ws1.Columns.Item(4).ColumnWidth = 7.330729166666667
$ws1.Columns.Item(5).ColumnWidth = 5.330729166666667

$ws3.Columns.Item(3).ColumnWidth = 21.666666666666668
$ws3.Columns.Item(4).ColumnWidth = 7.166666666666667
$ws3.Columns.Item(5).ColumnWidth = 5.830729166666667

# ---------------------------------------------------------------------------
# Selections / active-tab bookkeeping. The LAST sheet selected becomes the
# active tab, so Group B (the one that should end up tabSelected) goes last.
# ---------------------------------------------------------------------------
$ws1.Range("F14").Select()
$ws3.Range("E4:E18").Select()
$ws2.Range("F16").Select()
